$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 10990084.13351303
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10990333.33711062
